$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4963596
$ws.Range("B3").Value = 4963596
$ws.Range("B4").Value = 4963596
$ws.Range("B5").Value = 4963596
